$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update B9 and B10 values
$ws.Range("B9").Value = 2.37
$ws.Range("B10").Value = 2.93

# Apply the "underline" style (style index 1) to B11, matching existing
# cells like E8/E14 that use the same style (copy style from E8).
$ws.Range("E8").Copy()
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("E21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Move the active selection to B11, matching the diff's new selection
$ws.Range("B11").Select() | Out-Null
